# Insert a new daily price record as row 237 in the "Papa" (Feria Lagunitas
# de Puerto Montt) sheet. Inserting the row shifts every existing row from
# 237 downward by one position (old row 237 -> new row 238, ..., old row
# 290 -> new row 291), which matches the target diff (dimension grows from
# A1:R290 to A1:R291).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("237:237").Insert()

$ws.Cells.Item(237, 1).Value  = 4
$ws.Cells.Item(237, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(237, 3).Value  = 'Los Lagos'
$ws.Cells.Item(237, 4).Value  = 44543
$ws.Cells.Item(237, 5).Value  = 10
$ws.Cells.Item(237, 6).Value  = 100114001
$ws.Cells.Item(237, 7).Value  = 'Papa'
$ws.Cells.Item(237, 8).Value  = 'Pehuenche'
$ws.Cells.Item(237, 9).Value  = '1a nueva(o)'
$ws.Cells.Item(237, 10).Value = 150
$ws.Cells.Item(237, 11).Value = 11000
$ws.Cells.Item(237, 12).Value = 11000
$ws.Cells.Item(237, 13).Value = 11000
$ws.Cells.Item(237, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(237, 15).Value = 'Región de La Araucanía'
$ws.Cells.Item(237, 16).Value = 440
$ws.Cells.Item(237, 17).Value = 25
$ws.Cells.Item(237, 18).Value = 'Hortaliza'
